# Finalisation nettoyage : MAJ dictionnaire (ajout de "sale_id" et mise a
# jour des indicateurs "model") et extension de la table "columns".

$wb = $excel.ActiveWorkbook

# The worksheet whose tab is literally named "Sheet1" is the decoy/empty
# sheet in this workbook; the data-dictionary table actually lives on the
# tab named "variable" (it is the tabSelected / ActiveSheet one).
$ws = $wb.Worksheets.Item("variable")

# Target state (row -> variables, type, model) for rows 2..39. Row 1 (the
# header: variables/type/model) is unchanged.
$data = @(
    @("sale_id", "int64", "yes"),
    @("order_id", "object", "yes"),
    @("customer_id", "object", ""),
    @("order_status", "object", ""),
    @("order_purchase_timestamp", "object", "yes"),
    @("order_approved_at", "object", ""),
    @("order_delivered_carrier_date", "object", ""),
    @("order_delivered_customer_date", "object", "yes"),
    @("order_estimated_delivery_date", "object", ""),
    @("customer_unique_id", "object", "yes"),
    @("customer_zip_code_prefix", "int64", ""),
    @("customer_city", "object", ""),
    @("customer_state", "object", ""),
    @("order_item_id", "int64", ""),
    @("product_id", "object", ""),
    @("seller_id", "object", ""),
    @("shipping_limit_date", "object", "yes"),
    @("price", "float64", "yes"),
    @("freight_value", "float64", ""),
    @("payment_sequential", "int64", "yes"),
    @("payment_type", "object", "yes"),
    @("payment_installments", "int64", ""),
    @("payment_value", "float64", ""),
    @("review_id", "object", ""),
    @("review_score", "int64", "yes"),
    @("review_comment_title", "object", ""),
    @("review_comment_message", "object", ""),
    @("review_creation_date", "object", ""),
    @("review_answer_timestamp", "object", ""),
    @("product_category_name", "object", ""),
    @("product_name_lenght", "float64", ""),
    @("product_description_lenght", "float64", ""),
    @("product_photos_qty", "float64", "yes"),
    @("product_weight_g", "float64", ""),
    @("product_length_cm", "float64", ""),
    @("product_height_cm", "float64", ""),
    @("product_width_cm", "float64", ""),
    @("product_category_name_english", "object", "yes")
)

# Grow the "columns" table by one row (A1:C38 -> A1:C39), and prime the new
# row 39 by copying down the formatting of the last existing row so the new
# row matches the look of the rest of the table.
$tbl = $ws.ListObjects.Item("columns")
[void]$tbl.Resize($ws.Range("A1:C39"))
[void]$ws.Range("A38:C38").AutoFill($ws.Range("A38:C39"), 0)

# Re-write every data row (2..39) with the refreshed dictionary contents.
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Match column C's width to the rest of the (now 3-column) table.
[void]$ws.Columns("C").AutoFit()

# Restore the cursor to where the author left it after the refresh.
[void]$ws.Range("D6").Select()

# The hidden ExternalData_1 name (driving the Power Query refresh range)
# grows along with the table.
$wb.Names.Item("ExternalData_1").RefersTo = "=variable!`$A`$1:`$B`$39"
